# "added timepoints to harvester"
#
# The "Exposure conditions" sheet originally modeled a single timepoint
# (TP0) worth of harvester rows: 4 "chemical1"/BMD10 replicates, 4
# "CONTROL (SEE VEHICLE)" replicates, then 2 extraction blanks.
# This adds three more timepoints (TP1, TP2, TP3) - each repeating the
# same chemical/control replicate pattern - ahead of the (unchanged)
# extraction-blank rows, which move down to the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exposure conditions")

# Column layout from row 1: I=replicate, J=chemical name, K=dose, L=time point
$replicateCol = 9
$chemicalCol  = 10
$doseCol      = 11
$timeCol      = 12

$timePoints = "TP1", "TP2", "TP3"
$dataRow = 2

foreach ($tp in $timePoints) {
    for ($rep = 1; $rep -le 4; $rep++) {
        $ws.Cells.Item($dataRow, $replicateCol).Value = $rep
        $ws.Cells.Item($dataRow, $chemicalCol).Value  = "chemical1"
        $ws.Cells.Item($dataRow, $doseCol).Value      = "BMD10"
        $ws.Cells.Item($dataRow, $timeCol).Value      = $tp
        $dataRow++
    }
    for ($rep = 1; $rep -le 4; $rep++) {
        $ws.Cells.Item($dataRow, $replicateCol).Value = $rep
        $ws.Cells.Item($dataRow, $chemicalCol).Value  = "CONTROL (SEE VEHICLE)"
        $ws.Cells.Item($dataRow, $doseCol).Value      = 0
        $ws.Cells.Item($dataRow, $timeCol).Value      = $tp
        $dataRow++
    }
}

# The two original TP0 extraction-blank rows (previously rows 6-7) are
# unchanged in content, just pushed down to the bottom of the sheet.
for ($i = 0; $i -lt 2; $i++) {
    $ws.Cells.Item($dataRow, $replicateCol).Value = 0
    $ws.Cells.Item($dataRow, $chemicalCol).Value  = "EXTRACTION BLANK"
    # Force text "0" (not numeric 0) to match the original dose cell type.
    $ws.Cells.Item($dataRow, $doseCol).NumberFormat = "@"
    $ws.Cells.Item($dataRow, $doseCol).Value = "0"
    $ws.Cells.Item($dataRow, $timeCol).Value = "TP0"
    $dataRow++
}

Write-Host "Wrote rows 2..$($dataRow - 1)"
